$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.409.93"
$ws.Range("E2").Value = "  +1.07%  "

# Row 3
$ws.Range("D3").Value = "1.945.89"
$ws.Range("E3").Value = "  +2.50%  "

# Row 4
$ws.Range("E4").Value = "  +0.32%  "

# Row 5
$ws.Range("D5").Value = "'325.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7
$ws.Range("E7").Value = "  +0.72%  "

# Row 8
$ws.Range("D8").Value = "'0.3874"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "

# Row 9
$ws.Range("D9").Value = "'46.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10
$ws.Range("D10").Value = "'0.07839"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "

# Row 11
$ws.Range("D11").Value = "'0.9767"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.18%  "

# Row 12
$ws.Range("D12").Value = "'22.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.57%  "

# Row 13
$ws.Range("D13").Value = "1.927.51"
$ws.Range("E13").Value = "  +1.99%  "

# Row 14
$ws.Range("E14").Value = "  +0.67%  "

# Row 15
$ws.Range("D15").Value = "'5.761"
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.07055"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17
$ws.Range("D17").Value = "'86.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.06%  "

# Row 18
$ws.Range("D18").Value = "'1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "

# Row 19
$ws.Range("E19").Value = "  -0.86%  "

# Row 20
$ws.Range("D20").Value = "'17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.84%  "

# Row 21
$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "

# Row 22
$ws.Range("D22").Value = "29.444.97"
$ws.Range("E22").Value = "  +1.14%  "

# Row 23
$ws.Range("D23").Value = "'5.477"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.02%  "

# Row 24
$ws.Range("E24").Value = "  -0.35%  "

# Row 25
$ws.Range("D25").Value = "2.169.32"
$ws.Range("E25").Value = "  +2.42%  "

# Row 26
$ws.Range("D26").Value = "'2.099"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "

# Row 27
$ws.Range("D27").Value = "'157.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.78%  "

# Row 28
$ws.Range("D28").Value = "'19.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
$ws.Range("D29").Value = "'5.759"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.45%  "

# Row 30
$ws.Range("D30").Value = "'118.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "

# Row 32
$ws.Range("D32").Value = "'0.09364"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.41%  "

# Row 33
$ws.Range("D33").Value = "'0.8628"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.16%  "

# Row 34
$ws.Range("E34").Value = "  -0.80%  "

# Row 35
$ws.Range("E35").Value = "  -0.92%  "

# Row 36
$ws.Range("D36").Value = "'3.127"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "

# Row 37
$ws.Range("D37").Value = "'0.05771"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "

# Row 38
$ws.Range("D38").Value = "'1.155"
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "  +0.10%  "

# Row 40
$ws.Range("D40").Value = "'7.695"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "

# Row 41
$ws.Range("D41").Value = "'0.5671"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("E42").Value = "  -0.60%  "

# Row 43
$ws.Range("D43").Value = "'9.422"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "

# Row 44
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.733"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.23%  "

# Row 45
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "'0.000002785"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +32.49%  "

# Row 46
$ws.Range("D46").Value = "'0.5298"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.93%  "

# Row 47
$ws.Range("D47").Value = "'11.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.80%  "

# Row 48
$ws.Range("D48").Value = "'0.06868"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.95%  "

# Row 49
$ws.Range("D49").Value = "'2.087"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.77%  "

# Row 50
$ws.Range("D50").Value = "'1.819"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "

# Row 51
$ws.Range("D51").Value = "'111.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.10%  "
